$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 86
$ws1.Range("F5").Value = 304
$ws1.Range("F7").Value = 136
$ws1.Range("F11").Value = 40
$ws1.Range("F12").Value = 126
$ws1.Range("F13").Value = 2668
$ws1.Range("F15").Value = 35
$ws1.Range("F16").Value = 62
$ws1.Range("F17").Value = 25
$ws1.Range("F18").Value = 44
$ws1.Range("F19").Value = 550
$ws1.Range("F20").Value = 18
$ws1.Range("F21").Value = 630
$ws1.Range("F23").Value = 101
$ws1.Range("F24").Value = 52
$ws1.Range("F25").Value = 26
$ws1.Range("F26").Value = 57
$ws1.Range("F27").Value = 2241
$ws1.Range("F28").Value = 4797
$ws1.Range("F32").Value = 1244
$ws1.Range("F33").Value = 251
$ws1.Range("F34").Value = 2159
$ws1.Range("F35").Value = 571
$ws1.Range("F38").Value = 57
$ws1.Range("F39").Value = 141
$ws1.Range("F40").Value = 303
$ws1.Range("F41").Value = 445
$ws1.Range("F42").Value = 751
$ws1.Range("F43").Value = 15
$ws1.Range("F46").Value = 440

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 50

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 86
$ws4.Range("F5").Value = 304
$ws4.Range("F7").Value = 136
$ws4.Range("F11").Value = 40
$ws4.Range("F12").Value = 126
$ws4.Range("F13").Value = 2668
$ws4.Range("F15").Value = 35
$ws4.Range("F16").Value = 62
$ws4.Range("F17").Value = 50
$ws4.Range("F18").Value = 25
$ws4.Range("F19").Value = 44
$ws4.Range("F20").Value = 550
$ws4.Range("F21").Value = 18
$ws4.Range("F22").Value = 630
$ws4.Range("F24").Value = 101
$ws4.Range("F25").Value = 52
$ws4.Range("F26").Value = 26
$ws4.Range("F27").Value = 57
$ws4.Range("F28").Value = 2241
$ws4.Range("F29").Value = 4797
$ws4.Range("F33").Value = 1244
$ws4.Range("F34").Value = 251
$ws4.Range("F35").Value = 2159
$ws4.Range("F36").Value = 571
$ws4.Range("F39").Value = 57
$ws4.Range("F40").Value = 141
$ws4.Range("F41").Value = 303
$ws4.Range("F42").Value = 445
$ws4.Range("F43").Value = 751
$ws4.Range("F44").Value = 15
$ws4.Range("F47").Value = 440
